$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell that keeps the default body style (s=0) so we can
# re-apply it to cells whose style gets nudged by the text-prefix trick
# used below to force numeric-looking / empty strings to stay text.
$bodyStyle = $ws.Range("D2").Style

# --- Row 2 updates (text/string values) ---
$ws.Range("A2").Value = "01a9bd40-978a-46ab-85ed-9d603452834e"
$ws.Range("B2").Value = "HEAWYN"
$ws.Range("C2").Value = "Adidas"
$ws.Range("D2").Value = "Đức"
$ws.Range("E2").Value = "Da"
$ws.Range("F2").Value = "Giày chạy bộ"
$ws.Range("G2").Value = "Đế đinh"
$ws.Range("H2").Value = "white"

# Numeric-looking text values: a leading apostrophe forces Excel to keep
# them as text (matching t="s" in the target) instead of auto-converting
# to numbers; then restore the original (unstyled) look.
$ws.Range("I2").Value = "'550"
$ws.Range("I2").Style = $bodyStyle

$ws.Range("J2").Value = "'2400000"
$ws.Range("J2").Style = $bodyStyle

$ws.Range("K2").Value = "'550000"
$ws.Range("K2").Style = $bodyStyle

$ws.Range("L2").Value = "'12"
$ws.Range("L2").Style = $bodyStyle

# Boolean flags
$ws.Range("M2").Value = $false
$ws.Range("N2").Value = $false

# Image list cleared (no images for this error row) - keep it a text cell
$ws.Range("P2").Value = "'"
$ws.Range("P2").Style = $bodyStyle

$ws.Range("Q2").Value = "Sai định dạng hoặc để trống trường"

# Remove the second error row entirely
$ws.Rows(3).Delete()
